# Generate Report for Handoff
#
# The localization-status report records, per source file, the most recent
# handoff timestamp. A new handoff just completed for
# "5e7514f1-7f05-48b5-afce-eb7fdd1ff60a" (row 7 of every sheet), so its
# "Latest Handoff Date" / "Latest Handoff Datetime" cells are refreshed with
# the new handoff timestamps.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D ("Latest Handoff Date") for the row's source file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value = "2016-03-21 20:39:14"

# zh-cn sheet: column E ("Latest Handoff Datetime") for the same source file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value = "2016-03-21 20:39:10"

# de-de sheet: column E ("Latest Handoff Datetime") for the same source file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value = "2016-03-21 20:39:14"
